$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the anchor paragraph ("Google drive allowed a form...") that the
# three new bibliography paragraphs must be inserted after.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "Google drive allowed a form to be created that could be filled as a questionnaire by potential users to complete user testing.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the anchor paragraph for the new entry."
}
$anchorIndex = $anchor.Paragraphs.Item(1).Index

# Template paragraph that already carries the numbered-list formatting
# (pStyle ListParagraph, numPr ilvl=0/numId=2) used by the bibliography items.
$template = $d.Paragraphs.Item($anchorIndex).Previous.Previous.Previous

# ---------------------------------------------------------------------------
# Insert three new, empty paragraphs after the anchor paragraph. Each new
# paragraph inherits the anchor's own formatting (ListParagraph style,
# ind left=680, no direct numbering) so paragraphs 2 and 3 already come out
# correct; paragraph 1 gets list numbering applied explicitly afterwards.
# ---------------------------------------------------------------------------
$anchorPara = $d.Paragraphs.Item($anchorIndex)
$anchorPara.Range.InsertParagraphAfter()

$idxAdMob = $anchorIndex + 1
$paraAdMob = $d.Paragraphs.Item($idxAdMob)
$paraAdMob.Range.InsertParagraphAfter()

$idxLink = $anchorIndex + 2
$paraLink = $d.Paragraphs.Item($idxLink)
$paraLink.Range.InsertParagraphAfter()

$idxDesc = $anchorIndex + 3

# ---------------------------------------------------------------------------
# Paragraph 1: numbered bibliography entry for "Google AdMob Ads SDK".
# ---------------------------------------------------------------------------
$paraAdMob = $d.Paragraphs.Item($idxAdMob)
$r1 = $paraAdMob.Range
$ip1 = $d.Range($r1.Start, $r1.Start)
$ip1.Text = "Google "
$ip1.Collapse(0)
$ip1.Text = "AdMob Ads SDK"
$ip1.Collapse(0)
$ip1.Text = ". An SDK which provides APIs that allow mobile app developers to easily display adverts in their application to monetise the app."

# ---------------------------------------------------------------------------
# Paragraph 2: hyperlink to the AdMob SDK documentation.
# ---------------------------------------------------------------------------
$paraLink = $d.Paragraphs.Item($idxLink)
$r2 = $paraLink.Range
$ip2 = $d.Range($r2.Start, $r2.Start)
$ip2.Text = "https://developers.google.com/mobile-ads-sdk/"
$null = $d.Hyperlinks.Add($ip2, "https://developers.google.com/mobile-ads-sdk/", $null, $null, $ip2.Text)

# ---------------------------------------------------------------------------
# Paragraph 3: description of why the SDK was used.
# ---------------------------------------------------------------------------
$paraDesc = $d.Paragraphs.Item($idxDesc)
$r3 = $paraDesc.Range
$ip3 = $d.Range($r3.Start, $r3.Start)
$ip3.Text = "The SDK was found when researching how to easily monetise apps with advertising."

# ---------------------------------------------------------------------------
# Apply the numbered-list formatting (numPr ilvl=0 / numId=2) to paragraph 1
# only, matching the other "Google ..." bibliography entries in this list.
# ---------------------------------------------------------------------------
$paraAdMob = $d.Paragraphs.Item($idxAdMob)
$admobRange = $paraAdMob.Range
$admobRange.ListFormat.List = $template.Range.ListFormat.List
$admobRange.ListFormat.ListLevelNumber = $template.Range.ListFormat.ListLevelNumber

Write-Output "Inserted AdMob bibliography entry after paragraph $anchorIndex."
